# ifo GDP component analysis preprocessing:
# append the newest year-over-year forecast vector (forecast date 2025-11-25,
# base year 2025 / target year 2026) as a new row at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 continues the same table as rows 2:38 - clone the formatting of the
# last data row (date style in column A, etc.) down into the new row first.
$ws.Range("A38:E38").Copy($ws.Range("A39:E39"))

# Forecast-issue date (serial 45986 = 2025-11-25), y_0 / y_0_forecast /
# y_1 / y_1_forecast for the new vintage.
$ws.Range("A39").Value = 45986
$ws.Range("B39").Value = 2025
$ws.Range("C39").Value = 0.21940175159154141
$ws.Range("D39").Value = 2026
$ws.Range("E39").Value = -0.18831859814396609

# Leave the sheet selection on the freshly-added row, matching where the
# author was working when the file was saved.
$null = $ws.Range("B39:E39").Select()
